# Adding MPA test automation upload file
#
# The "Data" sheet of this mass-asset-change upload template holds sample
# rows. This updates the sample Main Asset Number / Asset Subnumber values
# and the sample Depreciation Key value used by the MPA test automation
# scenario.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Helper: write a numeric-looking identifier as TEXT (these asset / subnumber
# columns store their sample values as text, matching the template's other
# text fields) and then drop the cell back to the default "Normal" style so
# we don't leave a stray number-format behind -- only the cell's stored value
# should change.
function Set-TextValue([string]$cellRef, [string]$value) {
    $ws.Range($cellRef).Value = "'" + $value
    $ws.Range($cellRef).Style = "Normal"
}

# *Main Asset Number (12) sample values
Set-TextValue "C6"  "20000065"
Set-TextValue "C7"  "60000212"
Set-TextValue "C9"  "60000213"
Set-TextValue "C11" "60000214"

# *Asset Subnumber (4) sample values
Set-TextValue "D8"  "163"
Set-TextValue "D10" "164"

# Depreciation Key (4) sample value: SUL2 -> MANU for every data row
$ws.Cells.Replace("SUL2", "MANU")
